# Update the NATMI ligand-receptor output (Agt-Lrp2) with recomputed TPM-based
# statistics. A new sending cluster "ECs" is introduced, which shifts the
# previous two data rows down by one and changes every derived-specificity
# value across the table (because those specificities depend on the full set
# of clusters present).
#
# Final table (rows 2-4), columns A-T:
#   Row2: ECs   | Agt | Lrp2 | MuSCs | 2 | 0.6666666666666666 | 0.1214023333333333 | 0.364207 | 0.2856182748266287 | 0.2856182748266287 | 3 | 1 | 0.110028 | 0.330084 | 1 | 1 | 0.013357655932 | 0.120218903388 | 0.2856182748266287 | 0.2856182748266287
#   Row3: FAPs  | Agt | Lrp2 | MuSCs | 1 | 0.3333333333333333 | 0.174539            | 0.523617 | 0.4106307243130825 | 0.4106307243130825 | 3 | 1 | 0.110028 | 0.330084 | 1 | 1 | 0.019204177092 | 0.172837593828 | 0.4106307243130825 | 0.4106307243130825
#   Row4: MuSCs | Agt | Lrp2 | MuSCs | 3 | 1                  | 0.1291096666666667 | 0.387329 | 0.3037510008602889 | 0.3037510008602889 | 3 | 1 | 0.110028 | 0.330084 | 1 | 1 | 0.014205678404 | 0.127851105636 | 0.3037510008602889 | 0.3037510008602889

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 (new "ECs" sending-cluster row) ----
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Agt"
$ws.Range("C2").Value = "Lrp2"
$ws.Range("D2").Value = "MuSCs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1214023333333333
$ws.Range("H2").Value = 0.364207
$ws.Range("I2").Value = 0.2856182748266287
$ws.Range("J2").Value = 0.2856182748266287
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.110028
$ws.Range("N2").Value = 0.330084
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.013357655932
$ws.Range("R2").Value = 0.120218903388
$ws.Range("S2").Value = 0.2856182748266287
$ws.Range("T2").Value = 0.2856182748266287

# ---- Row 3 (was row 2: "FAPs" sending-cluster row, values recomputed) ----
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Agt"
$ws.Range("C3").Value = "Lrp2"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.174539
$ws.Range("H3").Value = 0.523617
$ws.Range("I3").Value = 0.4106307243130825
$ws.Range("J3").Value = 0.4106307243130825
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.110028
$ws.Range("N3").Value = 0.330084
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.019204177092
$ws.Range("R3").Value = 0.172837593828
$ws.Range("S3").Value = 0.4106307243130825
$ws.Range("T3").Value = 0.4106307243130825

# ---- Row 4 (was row 3: "MuSCs" sending-cluster row, values recomputed) ----
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Agt"
$ws.Range("C4").Value = "Lrp2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1291096666666667
$ws.Range("H4").Value = 0.387329
$ws.Range("I4").Value = 0.3037510008602889
$ws.Range("J4").Value = 0.3037510008602889
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.110028
$ws.Range("N4").Value = 0.330084
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.014205678404
$ws.Range("R4").Value = 0.127851105636
$ws.Range("S4").Value = 0.3037510008602889
$ws.Range("T4").Value = 0.3037510008602889

Write-Output "Applied Agt-Lrp2 TPM update: added ECs row, recomputed rows 2-4."
